# Mark additional project-plan checklist items as completed.
#
# The document tracks progress on a list of project requirements using
# strikethrough formatting (struck-through = done). This change marks five
# more items as completed:
#   - "Data needs to be in the compatible format in order to insert properly"
#   - "Clicking on Loan Display button should launch another page:"
#   - "Include a control that displays the total amount of loan amounts"
#   - "Use a LINQ statement that returns this sum"
#   - "Use an online loan calculator to verify the results of your application"
#
# The last three of those were previously highlighted in yellow; the yellow
# highlight is removed as part of marking them done.

$d = $word.ActiveDocument

function Mark-Done($text, [bool]$removeHighlight) {
    $para = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($candidate.Range.Text.Trim() -eq $text) {
            $para = $candidate
            break
        }
    }
    if ($para -eq $null) {
        throw "Paragraph not found: $text"
    }
    if ($removeHighlight) {
        $para.Range.HighlightColorIndex = 0
    }
    $para.Range.Font.StrikeThrough = 1
}

Mark-Done "Data needs to be in the compatible format in order to insert properly" $false
Mark-Done "Clicking on Loan Display button should launch another page:" $false
Mark-Done "Include a control that displays the total amount of loan amounts" $true
Mark-Done "Use a LINQ statement that returns this sum" $true
Mark-Done "Use an online loan calculator to verify the results of your application" $true

Write-Host "Done marking items as completed."
